$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final target values for rows 2..10 (Modelo, Comparaciones_Significativas, Proporcion_Sig, Mejor_N_Calib, ECRPS_Mejor)
$data = @(
    @("Block Bootstrapping", "4/10", 102.4, 20,  10.43818135404258),
    @("AREPD",               "1/10", 25.6,  20,  9.11424723548325),
    @("AV-MCPS",             "0/10", 0,     40,  2.997357133040603),
    @("DeepAR",               "0/10", 0,     40,  3.141501599903346),
    @("EnCQR-LSTM",          "0/10", 0,     40,  4.308451514263713),
    @("LSPM",                "0/10", 0,     200, 1.073439263101331),
    @("LSPMW",               "0/10", 0,     200, 1.074735978208519),
    @("MCPS",                "0/10", 0,     40,  2.890950949747069),
    @("Sieve Bootstrap",     "0/10", 0,     200, 0.5366137304080879)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
